$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (D2:H2)
$ws.Range("D2").Value = 0.606
$ws.Range("E2").Value = 1.158
$ws.Range("F2").Value = 1.119
$ws.Range("G2").Value = 0.187
$ws.Range("H2").Value = 1.645

# Row 3 updates (D3:H3)
$ws.Range("D3").Value = 0.606
$ws.Range("E3").Value = 1.176
$ws.Range("F3").Value = 1.151
$ws.Range("G3").Value = 0.249
$ws.Range("H3").Value = 1.848

# Row 4 updates (D4:H4)
$ws.Range("D4").Value = 0.618
$ws.Range("E4").Value = 1.141
$ws.Range("F4").Value = 1.328
$ws.Range("G4").Value = 0.316
$ws.Range("H4").Value = 1.949
